$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189, shifting existing rows 189-202 down to 190-203
$ws.Rows.Item(189).Insert()

# Populate the new row 189 with data
$ws.Cells.Item(189, 1).Value = 4
$ws.Cells.Item(189, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(189, 3).Value = "Los Lagos"
$ws.Cells.Item(189, 4).Value = 44585
$ws.Cells.Item(189, 5).Value = 10
$ws.Cells.Item(189, 6).Value = 100112044
$ws.Cells.Item(189, 7).Value = "Perejil"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 60
$ws.Cells.Item(189, 11).Value = 6000
$ws.Cells.Item(189, 12).Value = 6000
$ws.Cells.Item(189, 13).Value = 6000
$ws.Cells.Item(189, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(189, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(189, 16).Value = 3000
$ws.Cells.Item(189, 17).Value = 2
$ws.Cells.Item(189, 18).Value = "Hortaliza"
